# Auto-generated Excel COM-interop script
# Commit: Add data for 2023-11-04
# Updates cumulative YTD crime-count cells across the Citywide Totals,
# By Neighborhood summary, and per-neighborhood worksheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("C2").Value = 65
$ws.Range("G2").Value = 80
$ws.Range("J2").Value = 108
$ws.Range("E3").Value = 129
$ws.Range("C9").Value = 439
$ws.Range("D9").Value = 372
$ws.Range("E9").Value = 413
$ws.Range("F9").Value = 463
$ws.Range("G9").Value = 413
$ws.Range("H9").Value = 412
$ws.Range("I9").Value = 462
$ws.Range("B10").Value = 1198
$ws.Range("C10").Value = 1417
$ws.Range("D10").Value = 1615
$ws.Range("E10").Value = 1917
$ws.Range("F10").Value = 1920
$ws.Range("G10").Value = 846
$ws.Range("H10").Value = 534
$ws.Range("I10").Value = 771
$ws.Range("J10").Value = 642
$ws.Range("B11").Value = 1670
$ws.Range("C11").Value = 2007
$ws.Range("D11").Value = 2201
$ws.Range("E11").Value = 2536
$ws.Range("F11").Value = 2590
$ws.Range("G11").Value = 1473
$ws.Range("H11").Value = 1195
$ws.Range("I11").Value = 1551
$ws.Range("J11").Value = 1357

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("F8").Value = 122
$ws.Range("F9").Value = 182

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("B8").Value = 37
$ws.Range("I8").Value = 14
$ws.Range("J8").Value = 15
$ws.Range("B9").Value = 43
$ws.Range("I9").Value = 33
$ws.Range("J9").Value = 35

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("J2").Value = 6
$ws.Range("D8").Value = 42
$ws.Range("D9").Value = 73
$ws.Range("J9").Value = 54

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("C7").Value = 33
$ws.Range("B8").Value = 181
$ws.Range("C8").Value = 279
$ws.Range("D8").Value = 467
$ws.Range("E8").Value = 562
$ws.Range("F8").Value = 504
$ws.Range("G8").Value = 153
$ws.Range("J8").Value = 103
$ws.Range("B9").Value = 226
$ws.Range("C9").Value = 330
$ws.Range("D9").Value = 532
$ws.Range("E9").Value = 634
$ws.Range("F9").Value = 572
$ws.Range("G9").Value = 230
$ws.Range("J9").Value = 209

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("E3").Value = 6
$ws.Range("E5").Value = 4
$ws.Range("E7").Value = 19

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("I6").Value = 16
$ws.Range("J7").Value = 12
$ws.Range("I8").Value = 31
$ws.Range("J8").Value = 23

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("E5").Value = 16
$ws.Range("I5").Value = 9
$ws.Range("E6").Value = 24
$ws.Range("F6").Value = 31
$ws.Range("E7").Value = 42
$ws.Range("F7").Value = 55
$ws.Range("I7").Value = 24

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("G6").Value = 12
$ws.Range("H6").Value = 8
$ws.Range("D7").Value = 36
$ws.Range("D8").Value = 55
$ws.Range("G8").Value = 42
$ws.Range("H8").Value = 22

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("E5").Value = 19
$ws.Range("D7").Value = 21
$ws.Range("F7").Value = 21
$ws.Range("E8").Value = 97
$ws.Range("G8").Value = 84
$ws.Range("H8").Value = 95
$ws.Range("B19").Value = 43
$ws.Range("I19").Value = 33
$ws.Range("J19").Value = 35
$ws.Range("C22").Value = 6
$ws.Range("G27").Value = 16
$ws.Range("C28").Value = 122
$ws.Range("E28").Value = 81
$ws.Range("F28").Value = 112
$ws.Range("F32").Value = 182
$ws.Range("D36").Value = 73
$ws.Range("J36").Value = 54
$ws.Range("E41").Value = 24
$ws.Range("I41").Value = 14
$ws.Range("C43").Value = 15
$ws.Range("I47").Value = 46
$ws.Range("J49").Value = 10
$ws.Range("E50").Value = 42
$ws.Range("F50").Value = 55
$ws.Range("I50").Value = 24
$ws.Range("D52").Value = 37
$ws.Range("B53").Value = 226
$ws.Range("C53").Value = 330
$ws.Range("D53").Value = 532
$ws.Range("E53").Value = 634
$ws.Range("F53").Value = 572
$ws.Range("G53").Value = 230
$ws.Range("J53").Value = 209
$ws.Range("C56").Value = 13
$ws.Range("D65").Value = 55
$ws.Range("G65").Value = 42
$ws.Range("H65").Value = 22
$ws.Range("H68").Value = 7
$ws.Range("I70").Value = 31
$ws.Range("J70").Value = 23
$ws.Range("B76").Value = 52
$ws.Range("C77").Value = 64
$ws.Range("E77").Value = 67
$ws.Range("F77").Value = 58
$ws.Range("H77").Value = 56
$ws.Range("E78").Value = 44
$ws.Range("E80").Value = 32
$ws.Range("F80").Value = 33
$ws.Range("H80").Value = 14
$ws.Range("E86").Value = 15
$ws.Range("E87").Value = 33
$ws.Range("C92").Value = 22
$ws.Range("J92").Value = 25
$ws.Range("E93").Value = 7
$ws.Range("B99").Value = 1670
$ws.Range("C99").Value = 2007
$ws.Range("D99").Value = 2201
$ws.Range("E99").Value = 2536
$ws.Range("F99").Value = 2590
$ws.Range("G99").Value = 1473
$ws.Range("H99").Value = 1195
$ws.Range("I99").Value = 1551
$ws.Range("J99").Value = 1357

$ws = $wb.Worksheets.Item('Sheffield & DePaul')
$ws.Range("F5").Value = 6
$ws.Range("E6").Value = 28
$ws.Range("H6").Value = 6
$ws.Range("E7").Value = 32
$ws.Range("F7").Value = 33
$ws.Range("H7").Value = 14

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("E5").Value = 17
$ws.Range("I5").Value = 6
$ws.Range("E6").Value = 24
$ws.Range("I6").Value = 14

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("E8").Value = 23
$ws.Range("E9").Value = 33

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range("E4").Value = 4
$ws.Range("E6").Value = 44

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("C7").Value = 39
$ws.Range("F7").Value = 37
$ws.Range("E8").Value = 52
$ws.Range("C9").Value = 122
$ws.Range("E9").Value = 81
$ws.Range("F9").Value = 112

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("I7").Value = 24
$ws.Range("I8").Value = 46

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("B7").Value = 49
$ws.Range("B8").Value = 52

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("E6").Value = 9
$ws.Range("E7").Value = 15

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("C8").Value = 19
$ws.Range("J8").Value = 20
$ws.Range("C9").Value = 22
$ws.Range("J9").Value = 25

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("D7").Value = 30
$ws.Range("D8").Value = 37

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("G2").Value = 2
$ws.Range("G6").Value = 16

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("E4").Value = 2

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("E6").Value = 7

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("C8").Value = 15
$ws.Range("E9").Value = 44
$ws.Range("F9").Value = 42
$ws.Range("H9").Value = 28
$ws.Range("C10").Value = 64
$ws.Range("E10").Value = 67
$ws.Range("F10").Value = 58
$ws.Range("H10").Value = 56

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("J5").Value = 6
$ws.Range("J6").Value = 10

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("C6").Value = 4

$ws = $wb.Worksheets.Item('Clearing')
$ws.Range("C7").Value = 6

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("D5").Value = 6
$ws.Range("F6").Value = 12
$ws.Range("D7").Value = 21
$ws.Range("F7").Value = 21

$ws = $wb.Worksheets.Item('Wrigleyville')
$ws.Range("H6").Value = 5
$ws.Range("H7").Value = 7

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("E6").Value = 34
$ws.Range("G6").Value = 22
$ws.Range("H7").Value = 38
$ws.Range("E8").Value = 97
$ws.Range("G8").Value = 84
$ws.Range("H8").Value = 95

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("C6").Value = 11
$ws.Range("C7").Value = 15

$ws = $wb.Worksheets.Item('North Park')
$ws.Range("C2").Value = 2

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range("C6").Value = 13
